# Refresh the "cryptos" price/volume table with the latest scrape values.
# Note: for Price-column (D) values that look like plain decimal numbers
# (e.g. "311.48"), the cell is forced to Text format before assignment
# (then reset to the Normal style) so Excel stores the literal digit
# string instead of silently converting it to a numeric cell - matching
# how these price strings are stored as text throughout the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.096.19"
$ws.Range("E2").Value = "  -1.00%  "

$ws.Range("D3").Value = "1.824.56"
$ws.Range("E3").Value = "  -0.85%  "

$ws.Range("E4").Value = "  -0.41%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.15%  "

$ws.Range("E6").Value = "  -0.44%  "

$ws.Range("E7").Value = "  -1.62%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3643"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.60%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07300"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.23%  "

$ws.Range("E10").Value = "  -1.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.17"
$ws.Range("D11").Style = "Normal"

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.893.81"
$ws.Range("E12").Value = "  +2.58%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07626"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.343"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.56%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.477"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.007"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.73%  "

$ws.Range("E18").Value = "  -2.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.009"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.37%  "

$ws.Range("D20").Value = "27.311.59"
$ws.Range("E20").Value = "  -0.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.194"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.58"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").Value = "2.079.17"
$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.867"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.06%  "

$ws.Range("E28").Value = "  -3.66%  "

$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.093"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.56%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.65%  "

$ws.Range("E31").Value = "  -0.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.959"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7332"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.78%  "

$ws.Range("E34").Value = "  -2.30%  "

$ws.Range("E35").Value = "  -3.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.010"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.32%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.533"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.43%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05268"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.91%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.072"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.04%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.942"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.09%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01916"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.126"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5220"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.61%  "

$ws.Range("E44").Value = "  -2.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.272"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4866"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.49%  "

$ws.Range("E47").Value = "  -0.45%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.88%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.641"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.29%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06254"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.27%  "
